$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6's formatting into rows 7, 8 and 9 (keeps existing style
# indices intact instead of creating brand-new cellXfs entries).
$ws.Rows("6").Copy()
$ws.Rows("7").Insert(-4121)
$ws.Rows("6").Copy()
$ws.Rows("8").Insert(-4121)
$ws.Rows("6").Copy()
$ws.Rows("9").Insert(-4121)

# Columns A-D of the new rows should carry no explicit style (only column E
# keeps the date style copied from row 6).
$ws.Range("A7:D7").ClearFormats()
$ws.Range("A8:D8").ClearFormats()
$ws.Range("A9:D9").ClearFormats()

# Row 7: 3151. Special Array I
$ws.Range("A7").Value = 3151
$ws.Range("B7").Value = "Special Array I"
$ws.Range("C7").Value = "Easy"
$ws.Range("D7").Value = "Arrays,Two pointers"
$ws.Range("E7").Value = 45689

# Row 8: 242. Valid Anagram
$ws.Range("A8").Value = 242
$ws.Range("B8").Value = "Valid Anagram"
$ws.Range("C8").Value = "Easy"
$ws.Range("D8").Value = "Arrays,Anagram,Hash Table"
$ws.Range("E8").Value = 45689

# Row 9: 1752. Check if Array Is Sorted and Rotated
$ws.Range("A9").Value = 1752
$ws.Range("B9").Value = "Check if Array Is Sorted and Rotated"
$ws.Range("C9").Value = "Easy"
$ws.Range("D9").Value = "Arrays,Cycle,Concatenation"
$ws.Range("E9").Value = 45690

# Column D is now wider to fit the longer "Topics" text.
$ws.Range("D1").EntireColumn.ColumnWidth = 21.95

# Selection moves on to the next empty row, ready for the next entry.
$ws.Range("B10").Select()
